$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("issues")
$ws2 = $wb.Worksheets.Item("assignees")

# --- "issues" sheet: re-assign two rows from peterfpeterson to AndreiSavici ---
$ws1.Range("C5").Value = "AndreiSavici"
$ws1.Range("C9").Value = "AndreiSavici"

# Refresh the data-validation list so the two re-assigned cells can still see
# "AndreiSavici" in their dropdown (it now lives one row further down the
# assignees sheet, so the source range grows by one row for just these cells).
$ws1.Range("C5").Validation.Delete()
$ws1.Range("C5").Validation.Add(3, 1, 1, "=assignees!`$A`$4:`$A`$22")
$ws1.Range("C9").Validation.Delete()
$ws1.Range("C9").Validation.Add(3, 1, 1, "=assignees!`$A`$4:`$A`$22")

# --- "assignees" sheet: move AndreiSavici (row 26) up to row 22, shifting
#     jmborr/mdoucet/gvardany/JeanBilheux down by one row ---
$ws2.Range("A22").Value = "AndreiSavici"
$ws2.Range("B22").Value = "Andrei Savici"
$ws2.Range("A23").Value = "jmborr"
$ws2.Range("B23").Value = "Jose Borreguero"
$ws2.Range("A24").Value = "mdoucet"
$ws2.Range("B24").Value = "Mathieu Doucet"
$ws2.Range("A25").Value = "gvardany"
$ws2.Range("B25").Value = "Gagik Vardanyan"
$ws2.Range("A26").Value = "JeanBilheux"
$ws2.Range("B26").Value = "Jean Bilheux"

# --- restore on-screen selections to match where the author left the cursor ---
[void]$ws2.Activate()
$ws2.Rows("22:22").Select() | Out-Null

[void]$ws1.Activate()
$ws1.Range("B15").Select() | Out-Null
